$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the previous data row down to the new row 23 first. This gives the
# new cells the same plain (unstyled) text type as the rest of the sheet,
# instead of Excel auto-detecting types (e.g. turning a "2025-06-19"
# string into a real date value/format) when assigning fresh values.
$ws.Range("A22:F22").Copy()
$ws.Range("A23:F23").PasteSpecial()

# Now overwrite just the cells that actually differ from row 22.
# F23 ("2025-06-19") already matches what was pasted from F22, so it is
# left untouched to avoid Excel's automatic date detection.
$ws.Range("A23").Value = 47
$ws.Range("B23").Value = "minor changes from edit2"
$ws.Range("E23").Value = "edit2 to main"
